# Movie rating project.pptx - update author name and college line on Slide 1
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "TextBox 3" shape (holds "Presented By: ..." text) by name,
# falling back to the known index if the name ever changes.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "TextBox 3") {
        $shape = $s.Shapes.Item($i)
        break
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(3)
}

$tr = $shape.TextFrame.TextRange

# --- 1) "1.DHARANI S" -> "DHARANI.S" ---
$full = $tr.Text
$oldName = "1.DHARANI S"
$newName = "DHARANI.S"
$idx = $full.IndexOf($oldName)
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $oldName.Length)
    $sub.Text = $newName
}

# --- 2) Split/replace the college line, changing department from
#        "CIVIL ENGINEERING" to "ELECTRICAL AND ELECTRONICS ENGINEERING"
#        and dropping the stray leading space before "VIVEKANANDHA" ---
$full = $tr.Text
$oldLine = " VIVEKANANDHA COLLEGE OF TECHNOLOGY FOR WOMEN-CIVIL ENGINEERING"
$idx = $full.IndexOf($oldLine)
if ($idx -ge 0) {
    $whole = $tr.Characters($idx + 1, $oldLine.Length)
    $whole.Text = "VIVEKANANDHA COLLEGE OF TECHNOLOGY FOR WOMEN-ELECTRICAL AND ELECTRONICS ENGINEERING"
}

# Re-set each logical piece individually so the line is stored as three
# separate runs (as produced when the department name is retyped in place).
$full = $tr.Text
$part1 = "VIVEKANANDHA "
$part2 = "COLLEGE OF TECHNOLOGY FOR "
$part3 = "WOMEN-ELECTRICAL AND ELECTRONICS ENGINEERING"

$idx1 = $full.IndexOf($part1)
if ($idx1 -ge 0) {
    $r1 = $tr.Characters($idx1 + 1, $part1.Length)
    $r1.Text = $part1
}

$idx2 = $full.IndexOf($part2)
if ($idx2 -ge 0) {
    $r2 = $tr.Characters($idx2 + 1, $part2.Length)
    $r2.Text = $part2
}

$idx3 = $full.IndexOf($part3)
if ($idx3 -ge 0) {
    $r3 = $tr.Characters($idx3 + 1, $part3.Length)
    $r3.Text = $part3
}
